# "Generate Report for Handback" — refresh the handback-status report's
# timestamp columns after a new handback cycle ran for the
# ad8f1733-723b-48bf-9ac9-c27b028504c0 source file.
#
# Overview!G2  = "Latest HO Xliff Generate Date" for ad8f1733...
# zh-cn!H2     = "Correspond Handoff Datetime"   for ad8f1733... (zh-cn)
# zh-cn!K2     = "Correspond Handback DateTime"  for ad8f1733... (zh-cn)
# de-de!H2     = "Correspond Handoff Datetime"   for ad8f1733... (de-de)
# de-de!K2     = "Correspond Handback DateTime"  for ad8f1733... (de-de)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-23 10:49:48"

$wsZhCn.Range("H2").Value = "2016-08-23 10:49:44"
$wsZhCn.Range("K2").Value = "2016-08-23 10:50:00"

$wsDeDe.Range("H2").Value = "2016-08-23 10:49:48"
$wsDeDe.Range("K2").Value = "2016-08-23 10:50:21"
